$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-7 (row 1 is the header, left untouched).
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..J = ligand stats, K..P = receptor stats, Q..T = edge stats.
$rows = @(
    @("sCs",   "ECs", "FAPs", "Ntng1",  2, 0.6666666666666666, 0.005643,            0.016929,  0.01716510603349066, 0.01716510603349066, 1, 0.3333333333333333, 0.011401,           0.034203, 0.01020808032994874, 0.01020808032994874, 0.00006433584299999999, 0.00057902258699999993, 0.0001752227812619606, 0.0001752227812619606),
    @("sCs",   "ECs", "FAPs", "Lrrc4c", 2, 0.6666666666666666, 0.005643,            0.016929,  0.01716510603349066, 0.01716510603349066, 3, 1,                  1.105459333333333, 3.316378, 0.9897919196700513,  0.9897919196700512,  0.006238107018,         0.05614296316199999,   0.0169898832522287,    0.0169898832522287),
    @("Ntng1", "ECs", "FAPs", "Ntng1",  1, 0.3333333333333333, 0.03066833333333334, 0.092005,  0.09328817890078024, 0.09328817890078024, 1, 0.3333333333333333, 0.011401,           0.034203, 0.01020808032994874, 0.01020808032994874, 0.0003496496683333333, 0.003146847015,         0.0009522932240537943, 0.0009522932240537943),
    @("Ntng1", "ECs", "FAPs", "Lrrc4c", 1, 0.3333333333333333, 0.03066833333333334, 0.092005,  0.09328817890078024, 0.09328817890078024, 3, 1,                  1.105459333333333, 3.316378, 0.9897919196700513,  0.9897919196700512,  0.03390259532111112,   0.30512335789,         0.09233588567672646,   0.09233588567672645),
    @("Lrrc4c","ECs", "FAPs", "Ntng1",  3, 1,                  0.292437,            0.877311,  0.8895467150657291,  0.8895467150657291,  1, 0.3333333333333333, 0.011401,           0.034203, 0.01020808032994874, 0.01020808032994874, 0.003334074237,        0.030006668133,         0.00908056432463299,   0.00908056432463299),
    @("Lrrc4c","ECs", "FAPs", "Lrrc4c", 3, 1,                  0.292437,            0.877311,  0.8895467150657291,  0.8895467150657291,  3, 1,                  1.105459333333333, 3.316378, 0.9897919196700513,  0.9897919196700512,  0.323277211062,        2.909494899558,         0.8804661507410962,    0.8804661507410961)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
